$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helpers
# ---------------------------------------------------------------------------

function Get-ParaPlainText($p) {
    return $p.Range.Text.TrimEnd([char]13, [char]7)
}

# Build a single-run <w:p> OOXML fragment (ListParagraph style, numId 1) for
# use with Range.InsertXML. $ilvl is the numbering level ("0" or "1").
function Build-ParaXml($text, $ilvl) {
    $escaped = $text -replace '&', '&amp;' -replace '<', '&lt;' -replace '>', '&gt;'
    $spacePreserve = ""
    $trimmed = $text.Trim()
    if ($text -ne $trimmed) {
        $spacePreserve = ' xml:space="preserve"'
    }
    $runXml = '<w:r><w:t' + $spacePreserve + '>' + $escaped + '</w:t></w:r>'
    $fragment = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="' + $ilvl + '"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
        $runXml +
        '</w:p></w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
    return $fragment
}

# Replace the full contents of every paragraph whose (trimmed) text equals
# $oldText with a single clean run containing $newText (drops any proofErr
# markers and merges multiple runs into one, matching what Word does when a
# sentence no longer needs spell-check split markers).
# NOTE: keep call sites to ONE argument per line assigned from a plain
# variable - the PS interpreter here mis-parses calls that receive more than
# one parenthesized sub-expression as arguments, so build strings ahead of
# time into locals and pass the locals.
function Replace-AllParasWithText($oldText, $newText, $ilvl) {
    $hits = New-Object System.Collections.ArrayList
    foreach ($p in $d.Paragraphs) {
        $cur = Get-ParaPlainText $p
        if ($cur -eq $oldText) {
            [void]$hits.Add($p)
        }
    }
    if ($hits.Count -eq 0) {
        Write-Host "NOT FOUND: $oldText"
        return
    }
    $xml = Build-ParaXml $newText $ilvl
    foreach ($p in $hits) {
        $p.Range.InsertXML($xml)
    }
}

# ---------------------------------------------------------------------------
# 1) Strip proofErr spell-check wrapping / merge runs (text itself unchanged
#    except where noted) for the package-install sub-bullets.
# ---------------------------------------------------------------------------

Replace-AllParasWithText "django" "django" "1"
Replace-AllParasWithText "djangorestframework" "djangorestframework" "1"
Replace-AllParasWithText "djangorestframework-simplejwt" "djangorestframework-simplejwt" "1"
Replace-AllParasWithText "mysql-connector-python" "mysql-connector-python" "1"

# ---------------------------------------------------------------------------
# 2) "Makemigrations" occurs twice (ilvl 0); both lose the proofErr wrap.
# ---------------------------------------------------------------------------

Replace-AllParasWithText "Makemigrations" "Makemigrations" "0"

# ---------------------------------------------------------------------------
# 3) Insert a brand-new checklist item "Register rest_framework in
#    INSTALLED_APPS" right before "Add JWT settings to project settings".
# ---------------------------------------------------------------------------

foreach ($p in $d.Paragraphs) {
    $cur = Get-ParaPlainText $p
    if ($cur -eq "Add JWT settings to project settings") {
        $p.Range.InsertParagraphBefore()
        break
    }
}

foreach ($p in $d.Paragraphs) {
    $cur = Get-ParaPlainText $p
    if ($cur -eq "") {
        $nxt = $p.Next()
        if ($nxt -ne $null) {
            $nxtText = Get-ParaPlainText $nxt
            if ($nxtText -eq "Add JWT settings to project settings") {
                $newXml = Build-ParaXml "Register rest_framework in INSTALLED_APPS" "0"
                $p.Range.InsertXML($newXml)
                break
            }
        }
    }
}

# ---------------------------------------------------------------------------
# 4) Merge multi-run sentences into single clean runs (drops proofErr too).
# ---------------------------------------------------------------------------

Replace-AllParasWithText "Create RegistrationSerializer" "Create RegistrationSerializer" "0"
Replace-AllParasWithText "Add registration url path to authentication app" "Add registration url path to authentication app" "0"
Replace-AllParasWithText "Add JWT url paths to authentication app" "Add JWT url paths to authentication app" "0"

$apos = [char]8217

$registerUrls = "Register urls in project" + $apos + "s urls file"
Replace-AllParasWithText $registerUrls $registerUrls "0"

Replace-AllParasWithText "Add CarsList class" "Add CarsList class" "0"

$carsAppUrlsFile = "Create Car" + $apos + "s app urls file"
Replace-AllParasWithText $carsAppUrlsFile $carsAppUrlsFile "0"

Replace-AllParasWithText "Create path to CarsList" "Create path to CarsList" "0"
Replace-AllParasWithText "Create get_all_cars function with allow all" "Create get_all_cars function with allow all" "0"

$addToPathCars = "Add to path to Car" + $apos + "s urls"
Replace-AllParasWithText $addToPathCars $addToPathCars "0"

Replace-AllParasWithText "Create get_cars with IsAuthenticated" "Create get_cars with IsAuthenticated" "0"

$addPathCars = "Add path to Car" + $apos + "s urls"
Replace-AllParasWithText $addPathCars $addPathCars "0"

Replace-AllParasWithText "Create create_car with IsAuthenticated" "Create create_car with IsAuthenticated" "0"

# ---------------------------------------------------------------------------
# 5) Move the <w:lastRenderedPageBreak/> marker: it used to sit on the first
#    run of "Create get_cars with IsAuthenticated"; now it sits on the
#    "Test" item immediately before it.
# ---------------------------------------------------------------------------

$prevTest = $null
foreach ($p in $d.Paragraphs) {
    $t = Get-ParaPlainText $p
    if ($t -eq "Create get_cars with IsAuthenticated") {
        if ($prevTest -ne $null) {
            $pkgXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
                '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
                '<pkg:xmlData>' +
                '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
                '<w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
                '<w:r><w:lastRenderedPageBreak/><w:t>Test</w:t></w:r>' +
                '</w:p></w:body></w:document>' +
                '</pkg:xmlData></pkg:part></pkg:package>'
            $prevTest.Range.InsertXML($pkgXml)
            break
        }
    }
    if ($t -eq "Test") {
        $prevTest = $p
    } else {
        $prevTest = $null
    }
}

Write-Host "edit complete"
